# typo fix for vercel Limit
# F2 ("est. 10,000 month") should read "est. 10,000 visitors month"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F2").Value = "est. 10,000 visitors month"

# Match the style of the surrounding left-aligned plain cells (e.g. F3/F4),
# dropping the stray fill/border flags that were left over on F2.
$ws.Range("F2").HorizontalAlignment = -4131  # xlLeft

# Leave the selection on the cell that was actually edited.
$ws.Range("F2").Select()
